# Append a hard-coded 1x3 table (styled "ColorfulList") to the end of the
# document body, right after the last paragraph ("first item in unordered
# list") and before the final section break.

$d = $word.ActiveDocument

$tableXml = @"
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="ColorfulList"/>
    <w:tblW w:type="auto" w:w="0"/>
    <w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="2880"/>
    <w:gridCol w:w="2880"/>
    <w:gridCol w:w="2880"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="2880"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>R_7pKN6mN761yWDKP</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="2880"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>5</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:type="dxa" w:w="2880"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>5</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
"@

$end = $d.Content.End
$insertRange = $d.Range($end, $end)
$insertRange.InsertXML($tableXml)
